# Generate Report for Handoff
# Adds two new "Ready for handoff" files (1404531e-... and 62d912f1-...)
# as new rows in the Overview, zh-cn and de-de tables, just before the
# existing d33b969f-... row (which shifts down to the bottom).

$wb = $excel.ActiveWorkbook

function Set-Text($ws, $addr, $text) {
    # Prefix with an apostrophe so Excel always stores a literal text value
    # (prevents "True"/"False" becoming booleans, etc.)
    $ws.Range($addr).Value = "'" + $text
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsO = $wb.Worksheets.Item("Overview")
$loO = $wsO.ListObjects.Item(1)
$loO.Resize($wsO.Range("A1:G5"))

Set-Text $wsO "A3" "1404531e-d15f-4bff-964a-b44737b59e94.md"
Set-Text $wsO "B3" "e2e\1404531e-d15f-4bff-964a-b44737b59e94.md"
Set-Text $wsO "C3" ".md"
Set-Text $wsO "D3" ""
Set-Text $wsO "E3" "Ready for handoff"
Set-Text $wsO "F3" "Ready for handoff"
Set-Text $wsO "G3" "2016-08-19 14:47:26"

Set-Text $wsO "A4" "62d912f1-9438-4068-b5d7-e343c97ec4a0.md"
Set-Text $wsO "B4" "e2e\62d912f1-9438-4068-b5d7-e343c97ec4a0.md"
Set-Text $wsO "C4" ".md"
Set-Text $wsO "D4" ""
Set-Text $wsO "E4" "Ready for handoff"
Set-Text $wsO "F4" "Ready for handoff"
Set-Text $wsO "G4" "2016-08-19 14:47:26"

Set-Text $wsO "A5" "d33b969f-376c-43ea-91f1-d200e39205d6.md"
Set-Text $wsO "B5" "e2e\d33b969f-376c-43ea-91f1-d200e39205d6.md"
Set-Text $wsO "C5" ".md"
Set-Text $wsO "D5" ""
Set-Text $wsO "E5" "Ready for handoff"
Set-Text $wsO "F5" "Ready for handoff"
Set-Text $wsO "G5" "2016-08-19 14:46:04"

# rebuild hyperlinks (column B) in final row order
$wsO.Range("A1").Hyperlinks.Delete()
$wsO.Hyperlinks.Add($wsO.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ad33de6bcb7cfc799d1296b9e9238e6377dce55b/e2e/d9bdb7b6-6cb9-46a2-bc56-2db6ac5b2633.md", "", "", "e2e\d9bdb7b6-6cb9-46a2-bc56-2db6ac5b2633.md") | Out-Null
$wsO.Hyperlinks.Add($wsO.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9c9c979942d5a660dff63774e77e9d995c830380/e2e/1404531e-d15f-4bff-964a-b44737b59e94.md", "", "", "e2e\1404531e-d15f-4bff-964a-b44737b59e94.md") | Out-Null
$wsO.Hyperlinks.Add($wsO.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9c9c979942d5a660dff63774e77e9d995c830380/e2e/62d912f1-9438-4068-b5d7-e343c97ec4a0.md", "", "", "e2e\62d912f1-9438-4068-b5d7-e343c97ec4a0.md") | Out-Null
$wsO.Hyperlinks.Add($wsO.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9c9c979942d5a660dff63774e77e9d995c830380/e2e/d33b969f-376c-43ea-91f1-d200e39205d6.md", "", "", "e2e\d33b969f-376c-43ea-91f1-d200e39205d6.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZ = $wb.Worksheets.Item("zh-cn")
$loZ = $wsZ.ListObjects.Item(1)
$loZ.Resize($wsZ.Range("A1:P5"))

Set-Text $wsZ "A3" "1404531e-d15f-4bff-964a-b44737b59e94.md"
Set-Text $wsZ "B3" ".md"
Set-Text $wsZ "C3" "Ready for handoff"
Set-Text $wsZ "D3" "e2e"
Set-Text $wsZ "E3" "ht"
Set-Text $wsZ "F3" "False"
Set-Text $wsZ "G3" "1404531e-d15f-4bff-964a-b44737b59e94.12a82e05e722907c7e2114ba2f1a7efefdd695aa.zh-cn.xlf"
Set-Text $wsZ "H3" "2016-08-19 14:47:21"
Set-Text $wsZ "I3" ""
Set-Text $wsZ "J3" ""
Set-Text $wsZ "K3" "0001-01-01 00:00:00"
Set-Text $wsZ "L3" ""
Set-Text $wsZ "M3" "True"
Set-Text $wsZ "N3" ""
Set-Text $wsZ "O3" "False"
Set-Text $wsZ "P3" ""

Set-Text $wsZ "A4" "62d912f1-9438-4068-b5d7-e343c97ec4a0.md"
Set-Text $wsZ "B4" ".md"
Set-Text $wsZ "C4" "Ready for handoff"
Set-Text $wsZ "D4" "e2e"
Set-Text $wsZ "E4" "ht"
Set-Text $wsZ "F4" "False"
Set-Text $wsZ "G4" "62d912f1-9438-4068-b5d7-e343c97ec4a0.3d4d22298b82ec8fd983ed6be84af1b9c9d66288.zh-cn.xlf"
Set-Text $wsZ "H4" "2016-08-19 14:47:21"
Set-Text $wsZ "I4" ""
Set-Text $wsZ "J4" ""
Set-Text $wsZ "K4" "0001-01-01 00:00:00"
Set-Text $wsZ "L4" ""
Set-Text $wsZ "M4" "True"
Set-Text $wsZ "N4" ""
Set-Text $wsZ "O4" "False"
Set-Text $wsZ "P4" ""

Set-Text $wsZ "A5" "d33b969f-376c-43ea-91f1-d200e39205d6.md"
Set-Text $wsZ "B5" ".md"
Set-Text $wsZ "C5" "Ready for handoff"
Set-Text $wsZ "D5" "e2e"
Set-Text $wsZ "E5" "ht"
Set-Text $wsZ "F5" "False"
Set-Text $wsZ "G5" "d33b969f-376c-43ea-91f1-d200e39205d6.c2dd6d4ce98a694f0f3388642aa99def201f5d70.zh-cn.xlf"
Set-Text $wsZ "H5" "2016-08-19 14:45:57"
Set-Text $wsZ "I5" ""
Set-Text $wsZ "J5" ""
Set-Text $wsZ "K5" "0001-01-01 00:00:00"
Set-Text $wsZ "L5" ""
Set-Text $wsZ "M5" "True"
Set-Text $wsZ "N5" ""
Set-Text $wsZ "O5" "False"
Set-Text $wsZ "P5" ""

$wsZ.Range("A1").Hyperlinks.Delete()
$wsZ.Hyperlinks.Add($wsZ.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ad33de6bcb7cfc799d1296b9e9238e6377dce55b/e2e/d9bdb7b6-6cb9-46a2-bc56-2db6ac5b2633.md", "", "", "d9bdb7b6-6cb9-46a2-bc56-2db6ac5b2633.md") | Out-Null
$wsZ.Hyperlinks.Add($wsZ.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/64f4a0f658fa9075e6e7c51d8a055e4684873be9/e2e/d9bdb7b6-6cb9-46a2-bc56-2db6ac5b2633.md", "", "", "d9bdb7b6-6cb9-46a2-bc56-2db6ac5b2633.md") | Out-Null
$wsZ.Hyperlinks.Add($wsZ.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9c9c979942d5a660dff63774e77e9d995c830380/e2e/1404531e-d15f-4bff-964a-b44737b59e94.md", "", "", "1404531e-d15f-4bff-964a-b44737b59e94.md") | Out-Null
$wsZ.Hyperlinks.Add($wsZ.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9c9c979942d5a660dff63774e77e9d995c830380/e2e/62d912f1-9438-4068-b5d7-e343c97ec4a0.md", "", "", "62d912f1-9438-4068-b5d7-e343c97ec4a0.md") | Out-Null
$wsZ.Hyperlinks.Add($wsZ.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9c9c979942d5a660dff63774e77e9d995c830380/e2e/d33b969f-376c-43ea-91f1-d200e39205d6.md", "", "", "d33b969f-376c-43ea-91f1-d200e39205d6.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsD = $wb.Worksheets.Item("de-de")
$loD = $wsD.ListObjects.Item(1)
$loD.Resize($wsD.Range("A1:P5"))

Set-Text $wsD "A3" "1404531e-d15f-4bff-964a-b44737b59e94.md"
Set-Text $wsD "B3" ".md"
Set-Text $wsD "C3" "Ready for handoff"
Set-Text $wsD "D3" "e2e"
Set-Text $wsD "E3" "ht"
Set-Text $wsD "F3" "False"
Set-Text $wsD "G3" "1404531e-d15f-4bff-964a-b44737b59e94.12a82e05e722907c7e2114ba2f1a7efefdd695aa.de-de.xlf"
Set-Text $wsD "H3" "2016-08-19 14:47:26"
Set-Text $wsD "I3" ""
Set-Text $wsD "J3" ""
Set-Text $wsD "K3" "0001-01-01 00:00:00"
Set-Text $wsD "L3" ""
Set-Text $wsD "M3" "True"
Set-Text $wsD "N3" ""
Set-Text $wsD "O3" "False"
Set-Text $wsD "P3" ""

Set-Text $wsD "A4" "62d912f1-9438-4068-b5d7-e343c97ec4a0.md"
Set-Text $wsD "B4" ".md"
Set-Text $wsD "C4" "Ready for handoff"
Set-Text $wsD "D4" "e2e"
Set-Text $wsD "E4" "ht"
Set-Text $wsD "F4" "False"
Set-Text $wsD "G4" "62d912f1-9438-4068-b5d7-e343c97ec4a0.3d4d22298b82ec8fd983ed6be84af1b9c9d66288.de-de.xlf"
Set-Text $wsD "H4" "2016-08-19 14:47:26"
Set-Text $wsD "I4" ""
Set-Text $wsD "J4" ""
Set-Text $wsD "K4" "0001-01-01 00:00:00"
Set-Text $wsD "L4" ""
Set-Text $wsD "M4" "True"
Set-Text $wsD "N4" ""
Set-Text $wsD "O4" "False"
Set-Text $wsD "P4" ""

Set-Text $wsD "A5" "d33b969f-376c-43ea-91f1-d200e39205d6.md"
Set-Text $wsD "B5" ".md"
Set-Text $wsD "C5" "Ready for handoff"
Set-Text $wsD "D5" "e2e"
Set-Text $wsD "E5" "ht"
Set-Text $wsD "F5" "False"
Set-Text $wsD "G5" "d33b969f-376c-43ea-91f1-d200e39205d6.c2dd6d4ce98a694f0f3388642aa99def201f5d70.de-de.xlf"
Set-Text $wsD "H5" "2016-08-19 14:46:04"
Set-Text $wsD "I5" ""
Set-Text $wsD "J5" ""
Set-Text $wsD "K5" "0001-01-01 00:00:00"
Set-Text $wsD "L5" ""
Set-Text $wsD "M5" "True"
Set-Text $wsD "N5" ""
Set-Text $wsD "O5" "False"
Set-Text $wsD "P5" ""

$wsD.Range("A1").Hyperlinks.Delete()
$wsD.Hyperlinks.Add($wsD.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ad33de6bcb7cfc799d1296b9e9238e6377dce55b/e2e/d9bdb7b6-6cb9-46a2-bc56-2db6ac5b2633.md", "", "", "d9bdb7b6-6cb9-46a2-bc56-2db6ac5b2633.md") | Out-Null
$wsD.Hyperlinks.Add($wsD.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/a4decfe75757e94199921324d3b7af5fb486f609/e2e/d9bdb7b6-6cb9-46a2-bc56-2db6ac5b2633.md", "", "", "d9bdb7b6-6cb9-46a2-bc56-2db6ac5b2633.md") | Out-Null
$wsD.Hyperlinks.Add($wsD.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9c9c979942d5a660dff63774e77e9d995c830380/e2e/1404531e-d15f-4bff-964a-b44737b59e94.md", "", "", "1404531e-d15f-4bff-964a-b44737b59e94.md") | Out-Null
$wsD.Hyperlinks.Add($wsD.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9c9c979942d5a660dff63774e77e9d995c830380/e2e/62d912f1-9438-4068-b5d7-e343c97ec4a0.md", "", "", "62d912f1-9438-4068-b5d7-e343c97ec4a0.md") | Out-Null
$wsD.Hyperlinks.Add($wsD.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9c9c979942d5a660dff63774e77e9d995c830380/e2e/d33b969f-376c-43ea-91f1-d200e39205d6.md", "", "", "d33b969f-376c-43ea-91f1-d200e39205d6.md") | Out-Null

Write-Output "Report for handoff generated."
